$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 8-15 (extr1..extr8): only C (from_bus), D (to_bus)
#     and occasionally E (in_service) change. Column A/B are untouched.

# row 8  - extr1
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# row 9  - extr2
$ws.Range("C9").Value = 16

# row 10 - extr3
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# row 11 - extr4
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# row 12 - extr5
$ws.Range("C12").Value = 10

# row 13 - extr6
$ws.Range("D13").Value = 8

# row 14 - extr7
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

# row 15 - extr8
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# --- Append two new rows (16, 17) for line7 / line8, matching the
#     formatting of the existing data rows (style on column A).

$ws.Range("A2").Copy($ws.Range("A16"))
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A2").Copy($ws.Range("A17"))
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true
